# Fix Training Data Issue
# The BF column ("Date") was storing the source filename-derived string
# "6-7-2013-14" (day-month confused with the season label) for every row.
# Correct it to the actual ISO-style game date "2014-06-07" for all data
# rows (BF2:BF31), leaving the header (BF1) and every other column intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateCol = 58   # column BF
$firstRow = 2
$lastRow = 31
$oldValue = "6-7-2013-14"
$newValue = "2014-06-07"

# Force the range to be treated as plain text first so the COM layer
# doesn't "helpfully" reinterpret the replacement string as a real date
# serial value (which would change both the stored type and the style).
$rng = $ws.Range($ws.Cells.Item($firstRow, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
$rng.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Restore the default cell style so no stray number-format/style survives
# on these cells (matches the untouched look of the rest of the sheet).
$rng.Style = "Normal"
